# Auto-generated Excel COM-interop script applying a scheduled market-data refresh
# to the Leve profit calculation columns (H-N) across all job sheets.
# Values are static (no formulas in this workbook) and come from an external price feed.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 3433.1667
$ws.Range("I41").Value = 6066.6665
$ws.Range("J41").Value = 799.6667
$ws.Range("K41").Value = 6066.6665
$ws.Range("L41").Value = 799.6667
$ws.Range("M41").Value = -5626.6665
$ws.Range("N41").Value = -1679.6667
$ws.Range("H43").Value = 6666.3335
$ws.Range("I43").Value = 1000
$ws.Range("K43").Value = 1000
$ws.Range("M43").Value = -931
$ws.Range("H62").Value = 2390.0833
$ws.Range("I62").Value = 2648.4443
$ws.Range("J62").Value = 1615
$ws.Range("K62").Value = 2648.4443
$ws.Range("L62").Value = 1615
$ws.Range("M62").Value = -2024.4443
$ws.Range("N62").Value = -2863
$ws.Range("H65").Value = 2390.0833
$ws.Range("I65").Value = 2648.4443
$ws.Range("J65").Value = 1615
$ws.Range("K65").Value = 13242.2215
$ws.Range("L65").Value = 8075
$ws.Range("M65").Value = -10122.2215
$ws.Range("N65").Value = -14315
$ws.Range("H87").Value = 54719.5
$ws.Range("J87").Value = 57959.332
$ws.Range("L87").Value = 57959.332
$ws.Range("N87").Value = -60455.332
$ws.Range("H90").Value = 54719.5
$ws.Range("J90").Value = 57959.332
$ws.Range("L90").Value = 173877.996
$ws.Range("N90").Value = -186357.996
$ws.Range("H116").Value = 4094.1667
$ws.Range("J116").Value = 4392.5557
$ws.Range("L116").Value = 4392.5557
$ws.Range("N116").Value = -11276.5557
$ws.Range("H125").Value = 10361.546
$ws.Range("I125").Value = 14568.143
$ws.Range("K125").Value = 131113.287
$ws.Range("M125").Value = -128653.287
$ws.Range("H132").Value = 1561.6842
$ws.Range("I132").Value = 1259.2069
$ws.Range("J132").Value = 2536.3333
$ws.Range("K132").Value = 3777.620699999999
$ws.Range("L132").Value = 7608.999899999999
$ws.Range("M132").Value = -1247.620699999999
$ws.Range("N132").Value = -12668.9999
$ws.Range("H141").Value = 4465.1
$ws.Range("I141").Value = 4206.375
$ws.Range("K141").Value = 12619.125
$ws.Range("M141").Value = -7439.125

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2182.1738
$ws.Range("I2").Value = 1499.55
$ws.Range("K2").Value = 1499.55
$ws.Range("M2").Value = -1386.55
$ws.Range("H32").Value = 8522.031999999999
$ws.Range("I32").Value = 6120.923
$ws.Range("K32").Value = 6120.923
$ws.Range("M32").Value = -5833.923
$ws.Range("H61").Value = 6349.5366
$ws.Range("I61").Value = 6717.057
$ws.Range("K61").Value = 6717.057
$ws.Range("M61").Value = -6505.057
$ws.Range("H74").Value = 2795.5854
$ws.Range("I74").Value = 2233.353
$ws.Range("J74").Value = 5526.4287
$ws.Range("K74").Value = 2233.353
$ws.Range("L74").Value = 5526.4287
$ws.Range("M74").Value = -1359.353
$ws.Range("N74").Value = -7274.4287
$ws.Range("H77").Value = 2795.5854
$ws.Range("I77").Value = 2233.353
$ws.Range("J77").Value = 5526.4287
$ws.Range("K77").Value = 11166.765
$ws.Range("L77").Value = 27632.1435
$ws.Range("M77").Value = -6798.764999999999
$ws.Range("N77").Value = -36368.14350000001
$ws.Range("H112").Value = 18900
$ws.Range("J112").Value = 18900
$ws.Range("L112").Value = 18900
$ws.Range("N112").Value = -21854
$ws.Range("H116").Value = 2182.1738
$ws.Range("I116").Value = 1499.55
$ws.Range("K116").Value = 1499.55
$ws.Range("M116").Value = 794.45
$ws.Range("H122").Value = 2978
$ws.Range("I122").Value = 2972.5
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 8917.5
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -6467.5
$ws.Range("N122").Value = -13900
$ws.Range("H132").Value = 2585.1667
$ws.Range("I132").Value = 2592.9092
$ws.Range("K132").Value = 7778.7276
$ws.Range("M132").Value = -5248.7276
$ws.Range("H135").Value = 81037
$ws.Range("J135").Value = 81037
$ws.Range("L135").Value = 81037
$ws.Range("N135").Value = -91177
$ws.Range("H136").Value = 6349.5366
$ws.Range("I136").Value = 6717.057
$ws.Range("K136").Value = 20151.171
$ws.Range("M136").Value = -17601.171
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2182.1738
$ws.Range("I3").Value = 1499.55
$ws.Range("K3").Value = 1499.55
$ws.Range("M3").Value = -1385.55
$ws.Range("H86").Value = 2545.5386
$ws.Range("I86").Value = 2119.4
$ws.Range("K86").Value = 2119.4
$ws.Range("M86").Value = -996.4000000000001
$ws.Range("H89").Value = 2545.5386
$ws.Range("I89").Value = 2119.4
$ws.Range("K89").Value = 10597
$ws.Range("M89").Value = -4981
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 500.94446
$ws.Range("I6").Value = 1246.4
$ws.Range("K6").Value = 1246.4
$ws.Range("M6").Value = -1133.4
$ws.Range("H18").Value = 63999.668
$ws.Range("J18").Value = 63999.668
$ws.Range("L18").Value = 63999.668
$ws.Range("N18").Value = -64459.668
$ws.Range("H115").Value = 45084.5
$ws.Range("J115").Value = 45084.5
$ws.Range("L115").Value = 45084.5
$ws.Range("N115").Value = -47434.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 38516100
$ws.Range("I4").Value = 37686800
$ws.Range("K4").Value = 113060400
$ws.Range("M4").Value = -113060288
$ws.Range("H23").Value = 97.545456
$ws.Range("I23").Value = 9.5
$ws.Range("J23").Value = 117.111115
$ws.Range("K23").Value = 28.5
$ws.Range("L23").Value = 351.333345
$ws.Range("M23").Value = 206.5
$ws.Range("N23").Value = -821.333345
$ws.Range("H121").Value = 59783.41
$ws.Range("I121").Value = 403.75
$ws.Range("J121").Value = 112565.336
$ws.Range("K121").Value = 1211.25
$ws.Range("L121").Value = 337696.008
$ws.Range("M121").Value = 98.75
$ws.Range("N121").Value = -340316.008
$ws.Range("H132").Value = 21740922
$ws.Range("I132").Value = 38462864
$ws.Range("J132").Value = 2400
$ws.Range("K132").Value = 346165776
$ws.Range("L132").Value = 21600
$ws.Range("M132").Value = -346163246
$ws.Range("N132").Value = -26660

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2518.2812
$ws.Range("I80").Value = 2404.0667
$ws.Range("K80").Value = 2404.0667
$ws.Range("M80").Value = -1406.0667
$ws.Range("H83").Value = 2518.2812
$ws.Range("I83").Value = 2404.0667
$ws.Range("K83").Value = 12020.3335
$ws.Range("M83").Value = -7028.333499999999
$ws.Range("H132").Value = 8111.846
$ws.Range("I132").Value = 7465.875
$ws.Range("J132").Value = 9145.4
$ws.Range("K132").Value = 22397.625
$ws.Range("L132").Value = 27436.2
$ws.Range("M132").Value = -19867.625
$ws.Range("N132").Value = -32496.2
$ws.Range("H134").Value = 133997.4
$ws.Range("J134").Value = 133997.4
$ws.Range("L134").Value = 401992.2
$ws.Range("N134").Value = -407062.2

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()
$ws.Range("H46").Value = 7328
$ws.Range("I46").Value = 6691.8096
$ws.Range("J46").Value = 10000
$ws.Range("K46").Value = 6691.8096
$ws.Range("L46").Value = 10000
$ws.Range("M46").Value = -6503.8096
$ws.Range("N46").Value = -10376
$ws.Range("H109").Value = 68000
$ws.Range("J109").Value = 68000
$ws.Range("L109").Value = 68000
$ws.Range("N109").Value = -70774
$ws.Range("H122").Value = 2943
$ws.Range("I122").Value = 2934.5
$ws.Range("J122").Value = 2965.6667
$ws.Range("K122").Value = 8803.5
$ws.Range("L122").Value = 8897.000100000001
$ws.Range("M122").Value = -6353.5
$ws.Range("N122").Value = -13797.0001
$ws.Range("H131").Value = 42998.332
$ws.Range("J131").Value = 42998.332
$ws.Range("L131").Value = 42998.332
$ws.Range("N131").Value = -53078.332
$ws.Range("H132").Value = 5782.121
$ws.Range("I132").Value = 5962.724
$ws.Range("K132").Value = 17888.172
$ws.Range("M132").Value = -15358.172

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2347.2327
$ws.Range("I122").Value = 1598.3055
$ws.Range("K122").Value = 4794.916499999999
$ws.Range("M122").Value = -2344.916499999999
$ws.Range("H132").Value = 5969.5557
$ws.Range("I132").Value = 5815.4165
$ws.Range("K132").Value = 17446.2495
$ws.Range("M132").Value = -14916.2495

